$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for "Cebollín" at
# "Vega Monumental Concepción". It becomes the newest row (row 60),
# pushing all the existing historical rows (old 60-77) down by one
# (new 61-78).
$ws.Rows.Item(60).Insert()

# Fill in the data for the newly inserted row 60.
$ws.Cells.Item(60, 1).Value  = 11
$ws.Cells.Item(60, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(60, 3).Value  = "Bíobío"
$ws.Cells.Item(60, 4).Value  = 44875
$ws.Cells.Item(60, 5).Value  = 8
$ws.Cells.Item(60, 6).Value  = 100112037
$ws.Cells.Item(60, 7).Value  = "Cebollín"
$ws.Cells.Item(60, 8).Value  = "Sin especificar"
$ws.Cells.Item(60, 9).Value  = "Primera"
$ws.Cells.Item(60, 10).Value = 80
$ws.Cells.Item(60, 11).Value = 2600
$ws.Cells.Item(60, 12).Value = 2700
$ws.Cells.Item(60, 13).Value = 2638
$ws.Cells.Item(60, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(60, 15).Value = "Región Metropolitana"
$ws.Cells.Item(60, 16).Value = 73
$ws.Cells.Item(60, 17).Value = 36
$ws.Cells.Item(60, 18).Value = "Hortaliza"
